$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 0.7861523333333333
$ws.Range("N2").Value = 2.358457
$ws.Range("O2").Value = 0.01668797875153133
$ws.Range("P2").Value = 0.01668797875153133
$ws.Range("Q2").Value = 158.3623237234662
$ws.Range("R2").Value = 1425.260913511196
$ws.Range("S2").Value = 0.008066134918526745
$ws.Range("T2").Value = 0.008066134918526743

$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.5736784050900728
$ws.Range("P3").Value = 0.5736784050900727
$ws.Range("Q3").Value = 5443.981362434284
$ws.Range("R3").Value = 48995.83226190856
$ws.Range("S3").Value = 0.2772874704719497
$ws.Range("T3").Value = 0.2772874704719497

$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 19.22475933333333
$ws.Range("N4").Value = 57.674278
$ws.Range("O4").Value = 0.4080918692916219
$ws.Range("P4").Value = 0.4080918692916219
$ws.Range("Q4").Value = 3872.630572935266
$ws.Range("R4").Value = 34853.67515641739
$ws.Range("S4").Value = 0.1972512145341716
$ws.Range("T4").Value = 0.1972512145341716

$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 0.07263
$ws.Range("N5").Value = 0.21789
$ws.Range("O5").Value = 0.00154174686677398
$ws.Range("P5").Value = 0.00154174686677398
$ws.Range("Q5").Value = 14.63056850988
$ws.Range("R5").Value = 131.67511658892
$ws.Range("S5").Value = 0.0007452033839912249
$ws.Range("T5").Value = 0.0007452033839912249

$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 0.7861523333333333
$ws.Range("N6").Value = 2.358457
$ws.Range("O6").Value = 0.01668797875153133
$ws.Range("P6").Value = 0.01668797875153133
$ws.Range("Q6").Value = 51.42801387321756
$ws.Range("R6").Value = 462.852124858958
$ws.Range("S6").Value = 0.002619469636083453
$ws.Range("T6").Value = 0.002619469636083452

$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("O7").Value = 0.5736784050900728
$ws.Range("P7").Value = 0.5736784050900727
$ws.Range("S7").Value = 0.09004884206676821
$ws.Range("T7").Value = 0.09004884206676821

$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("M8").Value = 19.22475933333333
$ws.Range("N8").Value = 57.674278
$ws.Range("O8").Value = 0.4080918692916219
$ws.Range("P8").Value = 0.4080918692916219
$ws.Range("Q8").Value = 1257.63309193757
$ws.Range("R8").Value = 11318.69782743813
$ws.Range("S8").Value = 0.06405714414298665
$ws.Range("T8").Value = 0.06405714414298665

$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 0.07263
$ws.Range("N9").Value = 0.21789
$ws.Range("O9").Value = 0.00154174686677398
$ws.Range("P9").Value = 0.00154174686677398
$ws.Range("Q9").Value = 4.75126319574
$ws.Range("R9").Value = 42.76136876166
$ws.Range("S9").Value = 0.0002420040895408411
$ws.Range("T9").Value = 0.0002420040895408411

$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 0.7861523333333333
$ws.Range("N10").Value = 2.358457
$ws.Range("O10").Value = 0.01668797875153133
$ws.Range("P10").Value = 0.01668797875153133
$ws.Range("Q10").Value = 47.51099152882978
$ws.Range("R10").Value = 427.5989237594681
$ws.Range("S10").Value = 0.00241995734069754
$ws.Range("T10").Value = 0.00241995734069754

$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("O11").Value = 0.5736784050900728
$ws.Range("P11").Value = 0.5736784050900727
$ws.Range("Q11").Value = 1633.273283141372
$ws.Range("R11").Value = 14699.45954827235
$ws.Range("S11").Value = 0.08319025858479037
$ws.Range("T11").Value = 0.08319025858479037

$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("M12").Value = 19.22475933333333
$ws.Range("N12").Value = 57.674278
$ws.Range("O12").Value = 0.4080918692916219
$ws.Range("P12").Value = 0.4080918692916219
$ws.Range("Q12").Value = 1161.845279981519
$ws.Range("R12").Value = 10456.60751983367
$ws.Range("S12").Value = 0.05917822220864345
$ws.Range("T12").Value = 0.05917822220864345

$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 0.07263
$ws.Range("N13").Value = 0.21789
$ws.Range("O13").Value = 0.00154174686677398
$ws.Range("P13").Value = 0.00154174686677398
$ws.Range("Q13").Value = 4.38938252604
$ws.Range("R13").Value = 39.50444273436
$ws.Range("S13").Value = 0.0002235718119790129
$ws.Range("T13").Value = 0.0002235718119790129

$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 0.7861523333333333
$ws.Range("N14").Value = 2.358457
$ws.Range("O14").Value = 0.01668797875153133
$ws.Range("P14").Value = 0.01668797875153133
$ws.Range("Q14").Value = 70.33354433418057
$ws.Range("R14").Value = 633.0018990076251
$ws.Range("S14").Value = 0.003582416856223592
$ws.Range("T14").Value = 0.003582416856223592

$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("O15").Value = 0.5736784050900728
$ws.Range("P15").Value = 0.5736784050900727
$ws.Range("Q15").Value = 2417.838381671125
$ws.Range("R15").Value = 21760.54543504013
$ws.Range("S15").Value = 0.1231518339665645
$ws.Range("T15").Value = 0.1231518339665645

$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("M16").Value = 19.22475933333333
$ws.Range("N16").Value = 57.674278
$ws.Range("O16").Value = 0.4080918692916219
$ws.Range("P16").Value = 0.4080918692916219
$ws.Range("Q16").Value = 1719.953507167973
$ws.Range("R16").Value = 15479.58156451175
$ws.Range("S16").Value = 0.0876052884058202
$ws.Range("T16").Value = 0.0876052884058202

$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 0.07263
$ws.Range("N17").Value = 0.21789
$ws.Range("O17").Value = 0.00154174686677398
$ws.Range("P17").Value = 0.00154174686677398
$ws.Range("Q17").Value = 6.497882291250001
$ws.Range("R17").Value = 58.48094062125001
$ws.Range("S17").Value = 0.0003309675812629014
$ws.Range("T17").Value = 0.0003309675812629014
